$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple D/E value updates (Price / Volume(1h)) ---
# NOTE: some "Price" values look like plain decimal numbers to Excel's
# auto-detection (e.g. "597.40"), so they must be forced to Text via
# NumberFormat="@" before assignment, otherwise Excel stores them as
# numeric values and strips trailing zeros / formatting.
$updates = @(
    @{ Row = 2;  D = "62.707.71";  DText = $false; E = "  +1.78%  " }
    @{ Row = 3;  D = "3.025.83";   DText = $false; E = "  +2.45%  " }
    @{ Row = 4;  E = "  -0.06%  " }
    @{ Row = 5;  D = "597.40";     DText = $true;  E = "  +2.65%  " }
    @{ Row = 6;  D = "150.57";     DText = $true;  E = "  +6.92%  " }
    @{ Row = 7;  D = "0.999";      DText = $true;  E = "  -0.12%  " }
    @{ Row = 8;  D = "3.024.96";   DText = $false; E = "  +2.52%  " }
    @{ Row = 9;  D = "0.520";      DText = $true;  E = "  -0.15%  " }
    @{ Row = 10; D = "6.41";       DText = $true;  E = "  +12.19%  " }
    @{ Row = 11; D = "0.151";      DText = $true;  E = "  +4.28%  " }
    @{ Row = 12; E = "  +0.74%  " }
    @{ Row = 13; D = "0.0000234";  DText = $true;  E = "  +3.90%  " }
    @{ Row = 14; D = "34.64";      DText = $true;  E = "  +2.47%  " }
    @{ Row = 15; E = "  +2.68%  " }
    @{ Row = 16; D = "3.522.11";   DText = $false; E = "  +2.41%  " }
    @{ Row = 19; D = "3.025.21";   DText = $false; E = "  +2.61%  " }
    @{ Row = 20; D = "449.52";     DText = $true;  E = "  +0.19%  " }
    @{ Row = 21; D = "14.22";      DText = $true;  E = "  +2.89%  " }
    @{ Row = 22; D = "0.692";      DText = $true;  E = "  +2.01%  " }
    @{ Row = 23; E = "  +2.64%  " }
    @{ Row = 24; D = "82.42";      DText = $true;  E = "  +1.73%  " }
    @{ Row = 27; D = "12.06";      DText = $true;  E = "  +0.07%  " }
    @{ Row = 28; D = "0.999";      DText = $true;  E = "  -0.10%  " }
    @{ Row = 29; D = "2.72";       DText = $true;  E = "  +4.24%  " }
    @{ Row = 30; E = "  -0.01%  " }
    @{ Row = 31; D = "7.23";       DText = $true;  E = "  +5.80%  " }
    @{ Row = 32; D = "2.16";       DText = $true;  E = "  +5.02%  " }
    @{ Row = 33; D = "27.62";      DText = $true;  E = "  +2.46%  " }
    @{ Row = 34; D = "0.110";      DText = $true;  E = "  +3.24%  " }
    @{ Row = 35; D = "0.0₃0860";   DText = $false; E = "  +11.26%  " }
    @{ Row = 36; E = "  +2.26%  " }
    @{ Row = 37; D = "5.86";       DText = $true;  E = "  +3.74%  " }
    @{ Row = 38; D = "3.05";       DText = $true;  E = "  +9.27%  " }
    @{ Row = 39; D = "2.08";       DText = $true;  E = "  +0.44%  " }
    @{ Row = 40; D = "50.19";      DText = $true;  E = "  +0.43%  " }
    @{ Row = 41; D = "9.03";       DText = $true;  E = "  -0.92%  " }
    @{ Row = 42; E = "  +2.87%  " }
    @{ Row = 43; D = "0.285";      DText = $true;  E = "  +9.07%  " }
    @{ Row = 44; D = "393.55";     DText = $true;  E = "  +1.33%  " }
    @{ Row = 45; D = "40.20";      DText = $true;  E = "  +9.43%  " }
    @{ Row = 46; E = "  +0.88%  " }
    @{ Row = 47; D = "2.745.89";   DText = $false; E = "  +1.29%  " }
    @{ Row = 48; D = "133.57";     DText = $true;  E = "  +2.97%  " }
    @{ Row = 50; E = "  +2.07%  " }
    @{ Row = 51; D = "0.108";      DText = $true;  E = "  +0.13%  " }
)

foreach ($entry in $updates) {
    $rowNum = $entry.Row
    if ($entry.ContainsKey("D")) {
        $cell = $ws.Range("D$rowNum")
        if ($entry.DText) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $entry.D
    }
    if ($entry.ContainsKey("E")) {
        $ws.Range("E$rowNum").Value = $entry.E
    }
}

# --- Row 17/18 swap: WrappedBTC <-> Polkadot (with new values) ---
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$d17 = $ws.Range("D17")
$d17.NumberFormat = "@"
$d17.Value = "7.03"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "62.602.96"
$ws.Range("E18").Value = "  +1.68%  "

# --- Row 25/26 swap: Fetch.AI <-> RenderToken (with new values) ---
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$d25 = $ws.Range("D25")
$d25.NumberFormat = "@"
$d25.Value = "10.93"
$ws.Range("E25").Value = "  +14.93%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$d26 = $ws.Range("D26")
$d26.NumberFormat = "@"
$d26.Value = "2.25"
$ws.Range("E26").Value = "  +5.18%  "
